# Update column H (8th column) values in Sheet1 to reflect the new
# velocity/size-based measurements (Turtlebot system data update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(1, 8).Value = 25.860997000000001
$ws.Cells.Item(2, 8).Value = 71.420990399999994
$ws.Cells.Item(3, 8).Value = 80.4113507
$ws.Cells.Item(4, 8).Value = 148.271277
$ws.Cells.Item(5, 8).Value = 77.787005300000004
$ws.Cells.Item(6, 8).Value = 318.72573240000003
$ws.Cells.Item(7, 8).Value = 348.08444429999997
$ws.Cells.Item(8, 8).Value = 81.882517500000006
$ws.Cells.Item(9, 8).Value = 89.736347600000002
$ws.Cells.Item(10, 8).Value = 199.23783420000001
$ws.Cells.Item(11, 8).Value = 244.98417230000001
$ws.Cells.Item(12, 8).Value = 245.85245330000001
$ws.Cells.Item(13, 8).Value = 256.17017679999998
$ws.Cells.Item(14, 8).Value = 10.634350400000001
$ws.Cells.Item(15, 8).Value = 110.98351220000001
$ws.Cells.Item(16, 8).Value = 89.6249349
$ws.Cells.Item(17, 8).Value = 249.88699460000001
